# Update the "Last Updated" date stamp from 21-Apr-20 to 23-Apr-20
# on every slide's first textbox, without disturbing any other
# runs/formatting in that text frame.

$p = $ppt.ActivePresentation

$oldDate = "21-Apr-20"
$newDate = "23-Apr-20"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)

    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)

        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            $fullText = $tr.Text

            $idx = $fullText.IndexOf($oldDate)
            while ($idx -ge 0) {
                $sub = $tr.Characters($idx + 1, $oldDate.Length)
                $sub.Text = $newDate

                $fullText = $tr.Text
                $idx = $fullText.IndexOf($oldDate)
            }
        }
    }
}
